$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row 1: "_old" -> "_FV2404", "_new" -> "_FV2410"
$headers = @(
  "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
  "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404",
  "diff",
  "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410",
  "Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the A1:U78 range into an Excel Table ("Table1") with autofilter
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U78"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
